$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap match data (columns F:V) between row 112 and row 115 ---
$cols = @(6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22)
foreach ($c in $cols) {
    $v1 = $ws.Cells.Item(112, $c).Value()
    $v2 = $ws.Cells.Item(115, $c).Value()
    $ws.Cells.Item(112, $c).Value = $v2
    $ws.Cells.Item(115, $c).Value = $v1
}

# --- 2) Swap match data (columns F:V) between row 116 and row 117 ---
foreach ($c in $cols) {
    $v1 = $ws.Cells.Item(116, $c).Value()
    $v2 = $ws.Cells.Item(117, $c).Value()
    $ws.Cells.Item(116, $c).Value = $v2
    $ws.Cells.Item(117, $c).Value = $v1
}

# --- 3) Append a new row 133 with the new match record ---
# Copy formatting/styles from the last existing row (132) down to the new row 133
$srcRow = $ws.Range("A132:V132")
$dstRow = $ws.Range("A133:V133")
$srcRow.Copy($dstRow)

$ws.Cells.Item(133, 1).Value = 132
$ws.Cells.Item(133, 2).Value = "romania"
$ws.Cells.Item(133, 3).Value = "liga-2"
$ws.Cells.Item(133, 4).Value = "2023-2024"
$ws.Cells.Item(133, 5).Value = 45254.5625
$ws.Cells.Item(133, 6).Value = "CSC Dumbravita"
$ws.Cells.Item(133, 7).Value = 0
$ws.Cells.Item(133, 8).Value = "Gloria Buzau"
$ws.Cells.Item(133, 9).Value = 1
$ws.Cells.Item(133, 10).Value = 4
$ws.Cells.Item(133, 11).Value = "23/11/2023 01:42"
$ws.Cells.Item(133, 12).Value = 3.68
$ws.Cells.Item(133, 13).Value = "24/11/2023 13:25"
$ws.Cells.Item(133, 14).Value = 3.32
$ws.Cells.Item(133, 15).Value = "23/11/2023 01:42"
$ws.Cells.Item(133, 16).Value = 3.41
$ws.Cells.Item(133, 17).Value = "24/11/2023 13:25"
$ws.Cells.Item(133, 18).Value = 1.82
$ws.Cells.Item(133, 19).Value = "23/11/2023 01:42"
$ws.Cells.Item(133, 20).Value = 2.03
$ws.Cells.Item(133, 21).Value = "24/11/2023 13:02"
$ws.Cells.Item(133, 22).Value = "https://www.betexplorer.com/football/romania/liga-2/csc-dumbravita-fc-buzau/4QjKL2sG/"
